$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the worksheet/tab
$ws.Name = "Coverage Analysis"

# Rewrite the header row (A1:O1)
$headers = @(
    "Grid Row",
    "Grid Column",
    "Grid Cell Area (cm²)",
    "Total Coverage (%)",
    "blueberry Coverage (%)",
    "blueberry Area (cm²)",
    "lingonberry Coverage (%)",
    "lingonberry Area (cm²)",
    "crowberry Coverage (%)",
    "crowberry Area (cm²)",
    "frame Coverage (%)",
    "frame Area (cm²)",
    "wood / stick Coverage (%)",
    "wood / stick Area (cm²)",
    "Total Detections"
)

for ($col = 1; $col -le $headers.Length; $col++) {
    $ws.Cells.Item(1, $col).Value = $headers[$col - 1]
}

# The original header row only spanned A1:I1; columns J1:O1 are brand new
# cells, so copy the existing header formatting (style "1") onto them.
$ws.Cells.Item(1, 1).Copy()
$ws.Range($ws.Cells.Item(1, 10), $ws.Cells.Item(1, 15)).PasteSpecial(-4122)

# Rewrite data rows 2-17: columns C..O become 156.25 (C) / 0 (D..O), A/B unchanged
for ($row = 2; $row -le 17; $row++) {
    $ws.Cells.Item($row, 3).Value = 156.25
    for ($col = 4; $col -le 15; $col++) {
        $ws.Cells.Item($row, $col).Value = 0
    }
}
